$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions (copy H1's formatting for consistent header style)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows: I and J columns
$values = @{
    2  = @(9, 9)
    3  = @(9, 9)
    4  = @(6, 7)
    5  = @(8, 9)
    6  = @(1, 2)
    7  = @(8, 9)
    8  = @(10, 10)
    9  = @(10, 10)
    10 = @(5, 5)
    11 = @(4, 5)
    12 = @(4, 5)
    13 = @(5, 6)
    14 = @(7, 7)
    15 = @(8, 8)
    16 = @(8, 8)
    17 = @(9, 9)
    18 = @(9, 9)
    19 = @(7, 7)
    20 = @(3, 3)
    21 = @(5, 5)
    22 = @(7, 7)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
